# Insert 3 new weekly rows of data before the current row 101.
# This shifts the existing rows 101-178 down to 104-181 (unchanged),
# and the 3 new rows become rows 101, 102 and 103 with fresh data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A101:R103").EntireRow.Insert()

# Common (template) values shared by every "Cebollín" row on this sheet.
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$catId     = 100112037
$categoria = "Cebollín"
$variedad  = "Sin especificar"
$clasif    = "Hortaliza"

# --- Row 101 ---
$ws.Cells.Item(101,1).Value  = 7
$ws.Cells.Item(101,2).Value  = $mercado
$ws.Cells.Item(101,3).Value  = $region
$ws.Cells.Item(101,4).Value  = 45126
$ws.Cells.Item(101,5).Value  = $codreg
$ws.Cells.Item(101,6).Value  = $catId
$ws.Cells.Item(101,7).Value  = $categoria
$ws.Cells.Item(101,8).Value  = $variedad
$ws.Cells.Item(101,9).Value  = "Primera"
$ws.Cells.Item(101,10).Value = 60
$ws.Cells.Item(101,11).Value = 7000
$ws.Cells.Item(101,12).Value = 7000
$ws.Cells.Item(101,13).Value = 7000
$ws.Cells.Item(101,14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(101,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(101,16).Value = 194
$ws.Cells.Item(101,17).Value = 36
$ws.Cells.Item(101,18).Value = $clasif

# --- Row 102 ---
$ws.Cells.Item(102,1).Value  = 7
$ws.Cells.Item(102,2).Value  = $mercado
$ws.Cells.Item(102,3).Value  = $region
$ws.Cells.Item(102,4).Value  = 45126
$ws.Cells.Item(102,5).Value  = $codreg
$ws.Cells.Item(102,6).Value  = $catId
$ws.Cells.Item(102,7).Value  = $categoria
$ws.Cells.Item(102,8).Value  = $variedad
$ws.Cells.Item(102,9).Value  = "Segunda"
$ws.Cells.Item(102,10).Value = 60
$ws.Cells.Item(102,11).Value = 6000
$ws.Cells.Item(102,12).Value = 6000
$ws.Cells.Item(102,13).Value = 6000
$ws.Cells.Item(102,14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(102,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(102,16).Value = 167
$ws.Cells.Item(102,17).Value = 36
$ws.Cells.Item(102,18).Value = $clasif

# --- Row 103 ---
$ws.Cells.Item(103,1).Value  = 7
$ws.Cells.Item(103,2).Value  = $mercado
$ws.Cells.Item(103,3).Value  = $region
$ws.Cells.Item(103,4).Value  = 45126
$ws.Cells.Item(103,5).Value  = $codreg
$ws.Cells.Item(103,6).Value  = $catId
$ws.Cells.Item(103,7).Value  = $categoria
$ws.Cells.Item(103,8).Value  = $variedad
$ws.Cells.Item(103,9).Value  = "Primera"
$ws.Cells.Item(103,10).Value = 150
$ws.Cells.Item(103,11).Value = 7000
$ws.Cells.Item(103,12).Value = 7000
$ws.Cells.Item(103,13).Value = 7000
$ws.Cells.Item(103,14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(103,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(103,16).Value = 194
$ws.Cells.Item(103,17).Value = 36
$ws.Cells.Item(103,18).Value = $clasif
